$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Transporte" budget line (row 7), mirroring the existing
# rows' look: bold 12pt label in column A, plain 12pt value in column B.
$ws.Range("A7").Value = "Transporte"
$ws.Range("B7").Value = "RS 49.780"

$ws.Range("A7").Font.Size = 12
$ws.Range("A7").Font.Bold = $true

$ws.Range("B7").Font.Size = 12

# Match the row height used by the rest of the table
$ws.Rows.Item(7).RowHeight = 15.75

# Reflect the selection left behind in the saved workbook
$ws.Range("C14").Select() | Out-Null
